# Auto-generated edit script: add translations to nl.wur.IWMPRAISE.xlsx
$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("main")
$wsTool = $wb.Worksheets.Item("IWMPRAISE_Tool")

# --- Sheet "main", row 2 (nl.wur.IWMPRAISE.1_0.name) new language columns ---
$wsMain.Range("C2").Value = "IWMPRAISE Integreret plantebeskyttelsesværktøj til håndtering af ukrudt"
$wsMain.Range("D2").Value = "IWMPRAISE Integriertes Unkrautmanagement-Tool"
$wsMain.Range("F2").Value = "IWMPRAISE Orodje za integrirano varstvo pred pleveli"
$wsMain.Range("H2").Value = "IWMPRAISE Integrerat verktyg för ogräshantering"
$wsMain.Range("J2").Value = "IWMPRAISE Piktžolių Kontrolės Įrankis"
$wsMain.Range("K2").Value = "IWMPRAISE Outil de gestion intégrée de la flore adventice"
$wsMain.Range("L2").Value = "Strumento di gestione integrata delle erbe infestanti IWMPRAISE"
$wsMain.Range("M2").Value = "IWMPRAISE Εργαλείο ολοκληρωμένης διαχείρισης ζιζανίων"

# --- Sheet "IWMPRAISE_Tool", row 2 (description) ---
$wsTool.Range("A2").Value = "nl.wur.IWMPRAISE.1_0.models.IWMPRAISE_Tool.description"
$wsTool.Range("B2").Value = "THE PEST: Perennial and annual weed infestations can lead to direct and indirect damage to a wide range of outdoor crops. `nTHE DECISION: Integrated weed management tactics may affect one or more axes in the weed life cycle. They may prevent the establishment of seedlings from the seedbank (axis 1), reduce the impact established weeds have on the crop (axis 2), or reduce the weed seed/bud return to the soil (axis 3).`nTHE IWMPRAISE TOOL: This tool supports users in identifying and understanding the IPM tools and tactics available to manage perennial and annual weeds in narrow row, broad row and perennial crops. Users can select/de-select options to identify the best approach for them. This tool is currently available in English only. `nSOURCE: The IWMPRAISE TOOL was developed as part of the EU funded Horizon 2020 IWMPRAISE (727321), 2017-2022. It is designed for use across Europe.  `nREFERENCE: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02`n "
$wsTool.Range("C2").Value = "SKADEVOLDER: Flerårige og enårige ukrudtsangreb kan føre til direkte og indirekte skader på en lang række afgrøder.`nBESLUTNING: Integreret ukrudtshåndtering kan påvirke et eller flere tidspunkter i ukrudtets livscyklus. De kan forhindre etablering af frøplanter fra frøbanken (akse 1), reducere den påvirkning, etableret ukrudt har på afgrøden (akse 2), eller reducere ukrudtsfrøet tilbage til jorden (akse 3). IWMPRAISE-VÆRKTØJET: Dette værktøj hjælper brugere med at identificere og forstå de IPM-værktøjer og -taktikker, der er tilgængelige til at håndtere flerårigt og enårigt ukrudt i smalle rækker, brede rækker og flerårige afgrøder. Brugere kan vælge/fravælge muligheder for at identificere den bedste tilgang til dem. Dette værktøj er i øjeblikket kun tilgængeligt på engelsk.`nKILDE: IWMPRAISE-VÆRKTØJET blev udviklet som en del af det EU-finansierede Horizon 2020 IWMPRAISE (727321), 2017-2022. Den er designet til brug i hele Europa.`nREFERENCE: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02"
$wsTool.Range("D2").Value = "SCHADERREGER: Der Befall mit mehrjährigen und einjährigen Unkräutern kann zu direkten und indirekten Schäden an einer Vielzahl von Freilandkulturen führen. `nENTSCHEIDUNG: Integrierte Unkrautbekämpfungsmaßnahmen können eine oder mehrere Achsen im Lebenszyklus von Unkraut beeinflussen. Sie können die Etablierung von Sämlingen aus der Samenbank verhindern (Achse 1), die Auswirkungen etablierter Unkräuter auf die Kulturpflanzen verringern (Achse 2) oder die Rückkehr von Unkrautsamen/Knospen in den Boden reduzieren (Achse 3).`nIWMPRAISE-TOOL: Dieses Tool unterstützt die Nutzer bei der Identifizierung und dem Verständnis der verfügbaren IPM-Instrumente und -Taktiken zur Bekämpfung von mehrjährigen und einjährigen Unkräutern in eng- und breitreihigen sowie mehrjährigen Kulturen. Die Benutzer können Optionen auswählen bzw. die Auswahl aufheben, um den für sie besten Ansatz zu ermitteln. Dieses Tool ist derzeit nur in englischer Sprache verfügbar. `nQUELLE: Das IWMPRAISE TOOL wurde im Rahmen des von der EU finanzierten Horizon 2020 IWMPRAISE (727321), 2017-2022, entwickelt. Es ist für den europaweiten Einsatz konzipiert.  `nREFERENZEN: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02`n "
$wsTool.Range("E2").Value = "DE AANTASTING: Meerjarige en eenjarige onkruiden kunnen leiden tot directe en indirecte schade aan een groot aantal gewassen in de vollegrond. `nDE BESLUIT: Tactieken voor geïntegreerde onkruidbestrijding kunnen een of meer assen in de levenscyclus van onkruid beïnvloeden. Zij kunnen de vestiging van zaailingen uit de aanwezige zaadvoorraad voorkomen (as 1), het effect van gevestigd onkruid op het gewas verminderen (as 2), of de terugkeer van onkruidzaad/knoppen naar de bodem verminderen (as 3).`nHET IWMPRAISE-TOOL: Dit hulpmiddel ondersteunt gebruikers bij het identificeren en begrijpen van de beschikbare IPM-instrumenten en -tactieken voor het beheer van meerjarige en eenjarige onkruiden op een smalle rij, een brede rij en in meerjarige gewassen. Gebruikers kunnen opties selecteren/de-selecteren om de voor hen beste aanpak te bepalen. Deze tool is momenteel alleen beschikbaar in het Engels. `nBRON: De IWMPRAISE TOOL is ontwikkeld als onderdeel van het door de EU gefinancierde Horizon 2020 IWMPRAISE (727321), 2017-2022. Het is ontworpen voor gebruik in heel Europa.  `nREFERENTIES: Kudsk et al. (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02."
$wsTool.Range("F2").Value = "ŠKODLJIV ORGANIZEM: Večletni in enoletni pleveli lahko povzročijo neposredno in posredno škodo na številnih kulturnih rastlinah.`nODLOČITEV: Taktike celostnega zatiranja plevelov lahko vplivajo na eno ali več osi v življenjskem ciklu plevelov. Preprečijo lahko razvoj semen iz semenske banke (os 1), zmanjšajo vpliv uveljavljenih plevelov na pridelek (os 2) ali zmanjšajo vračanje semen/plevelov v tla (os 3).`nOrodje IWMPRAISE pomaga uporabnikom pri prepoznavanju in razumevanju orodij in taktik IPM, ki so na voljo za zatiranje večletnih in enoletnih plevelov v ozkih ali širokih medvrstnih prostoorih in v večletnih posevkih. Uporabniki lahko izberejo kateri pristop je zanje najprimernejši. To orodje je trenutno na voljo samo v angleščini.`nVIR: IWMPRAISE TOOL je bilo razvito v okviru projekta IWMPRAISE (727321), ki ga financira EU v okviru programa Obzorje 2020, 2017-2022. Namenjeno je uporabi po vsej Evropi. `nREFERENCA: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02`n "
$wsTool.Range("G2").Value = "SKADEGJØREREN: Ettårige og flerårige ugras kan føre til direkte og indirekte skader på et bredt spekter av utendørs-kulturer.  AVGJØRELSEN: Integrert metoder for ugrashåndtering kan påvirke en eller flere deler (akser) av ugrasets livssyklus. De kan hindre etablering av frøplanter fra frøbanken (akse 1), redusere påvirkningen som etablert ugras har på avlingen (akse 2), eller redusere mengden ugras som kommer tilbake til jorda (akse 3). IWMPRAISE-VERKTØYET: Dette verktøyet støtter brukere i å identifisere og forstå IPV-verktøyene og -taktikkene som er tilgjengelige for å håndtere flerårig og ettårig ugress i smale rader, brede rader og flerårige kulturer. Brukere kan velge/fjerne valg for å finne den beste tilnærmingen for dem. Dette verktøyet er for øyeblikket kun tilgjengelig på engelsk.`nKILDE: IWMPRAISE TOOL ble utviklet som en del av det EU-finansierede Horizon 2020 IWMPRAISE (727321), 2017-2022. Det er designet for bruk over hele Europa. REFERANSE: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02`n"
$wsTool.Range("H2").Value = "SKADEGÖRAREN: Perenna och anuella ogräs kan leda till direkta och indirekta skador på ett brett spektrum av utomhusgrödor.  BESLUTET: Integrerad taktik för hantering av ogräsproblem kan påverka en eller flera axlar i ogräsets livscykel. De kan förhindra etablering av groddplantor från fröbanken (axel 1), minska den inverkan etablerade ogräs har på grödan (axel 2) eller minska tillförseln av ogräsets frö/knopp tillbaka till jorden (axel 3).  IWMPRAISE-VERKTYGET: Detta verktyg hjälper användare att identifiera och förstå de IPM-verktyg och tillvägagångssätt som finns tillgängliga för att hantera perenna och annuella ogräs i grödor odlade med små eller stora radavstånd och i fleråriga grödor. Användare kan välja/avmarkera alternativ för att identifiera den bästa metoden för dem. Det här verktyget är för närvarande endast tillgängligt på engelska.  KÄLLA: IWMPRAISE TOOL utvecklades som en del av det EU-finansierade Horizon 2020 IWMPRAISE-projektet (727321), 2017–2022. Den är utformad för användning i hela Europa.  REFERENS: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02"
$wsTool.Range("I2").Value = "KASVINTUHOOJA: Monivuotiset ja yksivuotiset rikkakasvit voivat aiheuttaa vahinkoa lukuisilla kasveille pelloilla ja puutarhoissa.`nPÄÄTÖS: Integroidut rikkakasvien hallintamenetelmät voivat vaikuttaa yhteen tai useampaan vaiheeseen rikkakasvien elinkaaressa. Ne voivat estää uusien taimien muodostumisen siemenpankista (1), vähentää rikkakasvien vaikutusta satokasveihin (2) tai vähentävät rikkakasvien siementen määrää (3).`nIWMPRAISE -työkalu: Tämä työkalu tukee käyttäjiä tunnistamaan ja ymmärtämään IPM-työkaluja ja taktiikoita monivuotisten ja yksivuotisten rikkakasvien hallitsemiseksi. Käyttäjät voivat valita vaihtoehtoja tunnistaakseen heille parhaan lähestymistavan. Tämä työkalu on tällä hetkellä saatavana vain englanniksi.`nLÄHDE: IWMPRAISE-työkalu kehitettiin osana EU:n rahoittamaa hanketta Horizon 2020 IWMPRAISE (727321), 2017-2022. Se on suunniteltu käytettäväksi kaikkialla Euroopassa.`nVIITE: Kudsk et al (2020) Outelows for Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02`n"
$wsTool.Range("J2").Value = "PIKTŽOLĖS: Daugiametės ir vienmetės piktžolės gali tiesiogiai ir netiesiogiai pakenkti įvairiems lauko augalams. SPRENDIMAS: Integruotos piktžolių kontrolės taktikos gali paveikti vieną ar daugiau piktžolių gyvavimo ciklo ašių. Jos gali užkirsti kelią daigų iš dirvožemio sėklų banko sudygimui (1 ašis), sumažinti esamų piktžolių poveikį pasėliams (2 ašis) arba sumažinti piktžolių sėklų ar pumpurų grįžimą į dirvožemį (3 ašis). IWMPRAISE ĮRANKIS: Šis įrankis padeda naudotojams nustatyti ir suprasti turimas IKOK priemones ir strategiją daugiametėms ir vienmetėms piktžolėms siauraeiliuose, plačiaeiliuose ir daugiamečiuose pasėliuose kontroliuoti. Naudotojai gali pasirinkti ir (arba) panaikinti parinkčių reikšmes, kad nustatytų jiems tinkamiausią metodą. Šiuo metu ši priemonė prieinama tik anglų kalba. ŠALTINIS: IWMPRAISE Įrankis buvo sukurtas įgyvendinant 2017-2022 m. ES finansuojamą programos `"Horizontas 2020`" projektą IWMPRAISE (727321). Jis skirtas naudoti visoje Europoje. NUORODOS: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02"
$wsTool.Range("K2").Value = "LE BIOAGRESSEUR : Les infestations de mauvaises herbes vivaces et annuelles peuvent causer des dommages directs et indirects à un large éventail de cultures. `nLA DÉCISION : Les tactiques de gestion intégrée des mauvaises herbes peuvent affecter un ou plusieurs processus du cycle de vie des mauvaises herbes. Elles peuvent empêcher l'établissement de plantules à partir de la banque de graines (axe 1), réduire l'impact des mauvaises herbes établies sur la culture (axe 2), ou réduire le retour au sol des graines/bourgeons de mauvaises herbes (axe 3).`nL'OUTIL IWMPRAISE : Cet outil aide les utilisateurs à identifier et à comprendre les outils et tactiques IPM disponibles pour gérer les mauvaises herbes vivaces et annuelles. Les utilisateurs peuvent sélectionner/désélectionner des options pour identifier la stratégie la plus adaptée. Cet outil est actuellement disponible en anglais seulement. `nSOURCE : L'outil IWMPRAISE TOOL a été développé dans le cadre du projet Horizon 2020 IWMPRAISE (727321), 2017-2022, financé par l'UE. Il est conçu pour être utilisé dans toute l'Europe.  `nRÉFÉRENCE : Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, `nhttps://doi.org/10.1564/v31_aug_02 "
$wsTool.Range("L2").Value = "L'INFESTANTE: le infestazioni di erbe infestanti perenni e annuali possono causare danni diretti e indiretti a un'ampia gamma di colture all'aperto. `nLA DECISIONE: Le tattiche di gestione integrata delle infestanti possono influire su uno o più assi del ciclo di vita delle infestanti. Possono impedire l'insediamento di piantine dalla banca dei semi (asse 1), ridurre l'impatto che le infestanti insediate hanno sulla coltura (asse 2) o ridurre il ritorno di semi/boccioli di infestanti nel terreno (asse 3).`nLo strumento IWMPRAISE: questo strumento aiuta gli utenti a identificare e comprendere gli strumenti e le tattiche IPM disponibili per gestire le infestanti perenni e annuali nelle colture a file strette, a file larghe e perenni. Gli utenti possono selezionare/deselezionare le opzioni per identificare l'approccio migliore per loro. Questo strumento è attualmente disponibile solo in inglese. `nFONTE: Lo strumento IWMPRAISE è stato sviluppato nell'ambito del progetto Horizon 2020 IWMPRAISE (727321), finanziato dall'UE nel periodo 2017-2022. È stato progettato per essere utilizzato in tutta Europa.  `nRIFERIMENTO: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02.`n "
$wsTool.Range("M2").Value = " Ο Εχθρός: Οι μολύνσεις από πολυετή και ετήσια ζιζάνια μπορούν να προκαλέσουν άμεσες και έμμεσες ζημιές σε ένα ευρύ φάσμα υπαίθριων καλλιεργειών. `nΗ ΑΠΟΦΑΣΗ: Οι τακτικές ολοκληρωμένης διαχείρισης ζιζανίων μπορούν να επηρεάσουν έναν ή περισσότερους άξονες στον κύκλο ζωής των ζιζανίων. Μπορούν να αποτρέψουν την εγκατάσταση σπορόφυτων από την τράπεζα σπόρων (άξονας 1), να μειώσουν τις επιπτώσεις που έχουν τα εγκατεστημένα ζιζάνια στην καλλιέργεια (άξονας 2) ή να μειώσουν την επιστροφή σπόρων/βλαστών ζιζανίων στο έδαφος (άξονας 3).`nΤΟ ΕΡΓΑΛΕΙΟ IWMPRAISE: Το εργαλείο αυτό υποστηρίζει τους χρήστες στον εντοπισμό και την κατανόηση των διαθέσιμων εργαλείων και τακτικών IPM για τη διαχείριση πολυετών και ετήσιων ζιζανίων σε στενόμακρες, πλατύμακρες και πολυετείς καλλιέργειες. Οι χρήστες μπορούν να επιλέξουν/αποεπιλέξουν επιλογές για να προσδιορίσουν την καλύτερη προσέγγιση για αυτούς. Αυτό το εργαλείο είναι προς το παρόν διαθέσιμο μόνο στην αγγλική γλώσσα. `nΠΗΓΗ: Το εργαλείο IWMPRAISE TOOL αναπτύχθηκε στο πλαίσιο του χρηματοδοτούμενου από την ΕΕ προγράμματος Horizon 2020 IWMPRAISE (727321), 2017-2022. Έχει σχεδιαστεί για χρήση σε ολόκληρη την Ευρώπη.  `nΑΝΑΦΟΡΑ: Kudsk et al (2020) Outlooks on Pest Management, 31, 152-159, https://doi.org/10.1564/v31_aug_02."

# --- Sheet "IWMPRAISE_Tool", row 3 (name) ---
$wsTool.Range("A3").Value = "nl.wur.IWMPRAISE.1_0.models.IWMPRAISE_Tool.name"
$wsTool.Range("B3").Value = "IWMPRAISE Weed Management Tool"
$wsTool.Range("C3").Value = "IWM PRAISE Ukrudtshåndteringsværktøj"
$wsTool.Range("D3").Value = "IWMPRAISE Unkrautmanagement-Tool"
$wsTool.Range("E3").Value = "IWMPRAISE Onkruid management tool"
$wsTool.Range("F3").Value = "IWMPRAISE Orodje za obvladovanje plevela"
$wsTool.Range("G3").Value = "Verktøy for ugrasbekjempelse"
$wsTool.Range("H3").Value = "Verktyg för ogräsbekämpning"
$wsTool.Range("I3").Value = "IWMPRAISE rikkakasvien hallintatyökalu"
$wsTool.Range("J3").Value = "IWMPRAISE Piktžolių Kontrolės Įrankis"
$wsTool.Range("K3").Value = "IWMPRAISE Outil de gestion des mauvaises herbes"
$wsTool.Range("L3").Value = "Strumento di gestione delle erbe infestanti IWMPRAISE"
$wsTool.Range("M3").Value = "IWMPRAISE Εργαλείο διαχείρισης ζιζανίων"

# --- Sheet "IWMPRAISE_Tool", row 4 (purpose) ---
$wsTool.Range("A4").Value = "nl.wur.IWMPRAISE.1_0.models.IWMPRAISE_Tool.purpose"
$wsTool.Range("B4").Value = "An online strategic tool for identifying suitable management practices for annual and perennial weeds in narrow row, broad row and perennial crops."
$wsTool.Range("C4").Value = "Et online strategisk værktøj til at identificere passende håndteringspraksis for en- og flerårigt ukrudt i smalle rækker, brede rækker og flerårige afgrøder."
$wsTool.Range("D4").Value = "Ein strategisches Online-Tool zur Ermittlung geeigneter Bewirtschaftungsmethoden für ein- und mehrjährige Unkräuter in schmalen und breiten Reihen sowie für mehrjährige Kulturen."
$wsTool.Range("E4").Value = "Een online strategische tool voor het identificeren van geschikte beheerspraktijken voor éénjarige en meerjarige onkruiden op een smalle rij,  een brede rij en voor meerjarige gewassen."
$wsTool.Range("F4").Value = "Spletno strateško orodje za prepoznavanje ustreznih praks upravljanja enoletnih in večletnih plevelov v posevkih z ozkimi in širokimi medvrstnimi prostori ter večletnih posevkih."
$wsTool.Range("G4").Value = "Et nettbasert strategisk verktøy for å finne passene plantevernstrategier mot ettårige og flerårige ugras i smale og brede rader, og i flerårige kulturer."
$wsTool.Range("H4").Value = "Ett nättbaserat strategiskt verktyg för hitta passande växtskyddsstrategier mot annuella och perenna ogräs i grödor odlade med små eller stora radavstånd och i fleråriga grödor"
$wsTool.Range("I4").Value = "Online-työkalu sopivien hallintamenetelmien tunnistamiseksi yksi- ja monivuotisille rikkakasveille eri tuotantokasveilla."
$wsTool.Range("J4").Value = "Internetinė strateginė priemonė, skirta išsirinkti tinkamas vienamečių ir daugiamečių piktžolių kontrolės praktikas siauraeilių, plačiaeilių ir daugiamečių augalų pasėliuose."
$wsTool.Range("K4").Value = "Un outil stratégique en ligne pour identifier les pratiques de gestion appropriées pour les mauvaises herbes annuelles et vivaces."
$wsTool.Range("L4").Value = "Uno strumento strategico online per identificare le pratiche di gestione adeguate per le infestanti annuali e perenni nelle colture a file strette, a file larghe e perenni."
$wsTool.Range("M4").Value = "Ένα διαδικτυακό στρατηγικό εργαλείο για τον εντοπισμό κατάλληλων πρακτικών διαχείρισης ετήσιων και πολυετών ζιζανίων σε στενές, πλατιές και πολυετείς καλλιέργειες."
